$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5, shifting existing rows 5-13 down to 6-14.
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new feedback entry.
$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = "Ele efetuou a limpeza dos locais enferrujados e ensinou como efetuar a limpeza do inox para não riscar .`nTirou foto dos produtos utilizados para analisar e ficou de dar um retorno.Serviço prestado de acordo."
$ws.Cells.Item(5, 3).NumberFormat = $ws.Cells.Item(6, 3).NumberFormat
$ws.Cells.Item(5, 3).Value = 45912.71314665509
$ws.Cells.Item(5, 4).Value = "Y2Y2NGI4MTktZmVjMi00YWEzLWE0NzctM2JiMDE0YzliZjU5OjU3MDE2"

# The auto row-height from inserting multi-line text shouldn't pin an
# explicit row height; re-run autofit so the row reverts to default sizing.
$ws.Rows.Item(5).EntireRow.AutoFit()
